$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "5000 UI/ML SOL INJ IV CX 1 EST PLAS X 1 FA VD TRANS X 5 ML (EMB HOSP) Ativo"
$ws.Range("E3").Value = "25 MG COM REV BL AL PLAS INC X 200 (EMB. HOSP.) Ativo"
$ws.Range("E4").Value = "5 MG COM CT BL AL PLAS TRANS X 20 Ativo"
$ws.Range("E6").Value = "2 MG/ML SOL GOT OR CT FR GOT PLAS OPC X 30 ML Ativo"
$ws.Range("E7").Value = "25 MG COM REV CT BL AL PLAS OPC X 20 Ativo"
$ws.Range("E8").Value = "30 MG COM CT BL AL PLAS PVC/PVDC TRANS X 30 Ativo"
$ws.Range("E9").Value = "25 MG COM REV CT BL AL PLAS PVDC TRANS X 20 Ativo"
$ws.Range("E10").Value = "3 MG COM REV CT BL AL PLAS TRANS X 10 Ativo"
$ws.Range("E11").Value = "1 MG COM REV CT BL AL PLAS TRANS X 10 Ativo"

$ws.Range("D13").Value = "Último registro encontrado: 102980249"
$ws.Range("E13").Value = "Não encontrado"
$ws.Range("F13").Value = "Pendente"

$ws.Range("E15").Value = "1,2 U/G POM DERM CT 10 BG AL X 30 G + ESP PLAS Ativo"
